# Refresh the auto-generated DSR test-profile identifiers on rows 2-4
# (DSR Name, Father/Spouse Name, Contact Number, Employee No) with a new
# batch of randomized values, as produced by the test-data generator.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Contact Number must stay text (leading zeros), so force a text format
# before writing the value - otherwise Excel would coerce it to a number.
$ws.Range("E2:E4").NumberFormat = "@"

# Row 2
$ws.Range("C2").Value = "AUTODSR_DD62B"
$ws.Range("D2").Value = "Father_06BC"
$ws.Range("E2").Value = "03675892100"
$ws.Range("G2").Value = "EMPE17D78"

# Row 3
$ws.Range("C3").Value = "AUTODSR_26B88"
$ws.Range("D3").Value = "Father_4835"
$ws.Range("E3").Value = "03678450100"
$ws.Range("G3").Value = "EMP3A0AC6"

# Row 4
$ws.Range("C4").Value = "AUTODSR_4D224"
$ws.Range("D4").Value = "Father_5343"
$ws.Range("E4").Value = "03680963400"
$ws.Range("G4").Value = "EMP218A1F"
